# Add a new row (A12) to the sheet with value "toto_12", reusing the
# same cell formatting as the other data rows (e.g. A11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row onto the new cell,
# then overwrite its value with the new text.
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "toto_12"

# Match the workbook's recorded selection/active cell after the edit.
$ws.Range("A12").Select()
